$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.909.82"
$ws.Range("E2").Value = "  +2.38%  "
$ws.Range("D3").Value = "2.217.35"
$ws.Range("E3").Value = "  +0.24%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "263.36"
$ws.Range("E5").Value = "  +2.35%  "
$ws.Range("D6").Value = "86.37"
$ws.Range("E6").Value = "  +12.34%  "
$ws.Range("D7").Value = "0.622"
$ws.Range("E7").Value = "  +2.04%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").Value = "0.607"
$ws.Range("E9").Value = "  +2.03%  "
$ws.Range("D10").Value = "46.40"
$ws.Range("E10").Value = "  +10.13%  "
$ws.Range("E11").Value = "  +1.84%  "
$ws.Range("D12").Value = "7.63"
$ws.Range("E12").Value = "  +9.51%  "
$ws.Range("D14").Value = "2.550.68"
$ws.Range("E14").Value = "  +0.29%  "
$ws.Range("D15").Value = "14.64"
$ws.Range("E15").Value = "  +0.95%  "
$ws.Range("D16").Value = "2.211.06"
$ws.Range("E16").Value = "  -0.34%  "
$ws.Range("D17").Value = "0.782"
$ws.Range("E17").Value = "  -0.33%  "
$ws.Range("D18").Value = "43.864.99"
$ws.Range("E18").Value = "  +2.39%  "
$ws.Range("E19").Value = "  +1.06%  "
$ws.Range("E20").Value = "  +0.43%  "
$ws.Range("E21").Value = "  -1.73%  "
$ws.Range("D22").Value = "2.40"
$ws.Range("E22").Value = "  +8.98%  "
$ws.Range("D23").Value = "232.32"
$ws.Range("E23").Value = "  +0.83%  "
$ws.Range("D24").Value = "9.07"
$ws.Range("E24").Value = "  -2.95%  "
$ws.Range("D26").Value = "10.82"
$ws.Range("E26").Value = "  +0.35%  "
$ws.Range("E27").Value = "  +5.15%  "
$ws.Range("D28").Value = "39.75"
$ws.Range("E28").Value = "  -5.74%  "
$ws.Range("E29").Value = "  +2.87%  "
$ws.Range("E30").Value = "  +1.73%  "
$ws.Range("D31").Value = "175.01"
$ws.Range("E31").Value = "  +0.65%  "
$ws.Range("B32").Value = "Hedera"
$ws.Range("C32").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D32").Value = "0.0888"
$ws.Range("E32").Value = "  +1.57%  "
$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D33").Value = "20.55"
$ws.Range("E33").Value = "  +1.10%  "
$ws.Range("D34").Value = "5.43"
$ws.Range("E34").Value = "  +4.01%  "
$ws.Range("E35").Value = "  +0.94%  "
$ws.Range("D36").Value = "0.112"
$ws.Range("E36").Value = "  +4.04%  "
$ws.Range("D37").Value = "0.0361"
$ws.Range("E37").Value = "  +0.20%  "
$ws.Range("E38").Value = "  +3.32%  "
$ws.Range("D39").Value = "3.28"
$ws.Range("E39").Value = "  +15.26%  "
$ws.Range("D40").Value = "12.44"
$ws.Range("E40").Value = "  -3.84%  "
$ws.Range("D41").Value = "64.71"
$ws.Range("E41").Value = "  +7.70%  "
$ws.Range("E42").Value = "  -0.77%  "
$ws.Range("D43").Value = "5.55"
$ws.Range("E43").Value = "  +4.63%  "
$ws.Range("D44").Value = "0.205"
$ws.Range("E44").Value = "  +2.64%  "
$ws.Range("D45").Value = "101.27"
$ws.Range("E45").Value = "  -1.44%  "
$ws.Range("B46").Value = "Cronos"
$ws.Range("C46").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D46").Value = "0.0985"
$ws.Range("E46").Value = "  +0.79%  "
$ws.Range("B47").Value = "FraxShare"
$ws.Range("C47").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D47").Value = "8.35"
$ws.Range("E47").Value = "  +0.21%  "
$ws.Range("E48").Value = "  +1.10%  "
$ws.Range("E49").Value = "  +4.56%  "
$ws.Range("B50").Value = "Stacks"
$ws.Range("C50").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D50").Value = "1.56"
$ws.Range("E50").Value = "  +8.68%  "
$ws.Range("B51").Value = "WOONetwork"
$ws.Range("C51").Value = "https://coinranking.com/coin/k-J3YwacF+woonetwork-woo"
$ws.Range("D51").Value = "0.448"
$ws.Range("E51").Value = "  -2.82%  "
